# Adds summary statistics rows (Average loss, Correlation to e100,
# Correlation to v100, en/vn correlation) below the existing per-run data
# table (rows 4-102, columns B-CY) on the active sheet, then moves the
# selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- helpers -------------------------------------------------------------

function Get-ColLetter($sheet, $colIndex) {
    $addr = $sheet.Cells.Item(1, $colIndex).Address($false, $false)
    return ($addr -replace '\d+$', '')
}

$firstDataCol = 2            # column B
$lastDataCol  = 103           # column CY
$firstDataRow = 4
$lastDataRow  = 102

# ---- row 104: "Average loss" ---------------------------------------------

$ws.Range("A104").Value = "Average loss"

for ($col = $firstDataCol; $col -le $lastDataCol; $col++) {
    $letter = Get-ColLetter $ws $col
    $formula = "=AVERAGE(" + $letter + $firstDataRow + ":" + $letter + $lastDataRow + ")"
    $ws.Cells.Item(104, $col).Formula = $formula
}

# ---- row 106: "Correlation to e100" (vs column AZ, $-anchored) ----------

$ws.Range("A106").Value = "Correlation to e100"
$anchorCol106 = Get-ColLetter $ws 52   # AZ

for ($col = $firstDataCol; $col -le $lastDataCol; $col++) {
    $letter = Get-ColLetter $ws $col
    $formula = "=CORREL(" + $letter + $firstDataRow + ":" + $letter + $lastDataRow + ",`$" + $anchorCol106 + $firstDataRow + ":`$" + $anchorCol106 + $lastDataRow + ")"
    $ws.Cells.Item(106, $col).Formula = $formula
}

# ---- row 107: "Correlation to v100" (vs column CY, $-anchored) ----------

$ws.Range("A107").Value = "Correlation to v100"
$anchorCol107 = Get-ColLetter $ws 103  # CY

for ($col = $firstDataCol; $col -le $lastDataCol; $col++) {
    $letter = Get-ColLetter $ws $col
    $formula = "=CORREL(" + $letter + $firstDataRow + ":" + $letter + $lastDataRow + ",`$" + $anchorCol107 + $firstDataRow + ":`$" + $anchorCol107 + $lastDataRow + ")"
    $ws.Cells.Item(107, $col).Formula = $formula
}

# ---- row 109: "en/vn correlation" (column X vs column X+51) --------------

$ws.Range("A109").Value = "en/vn correlation"
$pairOffset = 51

for ($col = $firstDataCol; $col -le 52; $col++) {
    $letter = Get-ColLetter $ws $col
    $pairCol = $col + $pairOffset
    $pairLetter = Get-ColLetter $ws $pairCol
    $formula = "=CORREL(" + $letter + $firstDataRow + ":" + $letter + $lastDataRow + ", " + $pairLetter + $firstDataRow + ":" + $pairLetter + $lastDataRow + ")"
    $ws.Cells.Item(109, $col).Formula = $formula
}

# ---- move selection to where the author left it --------------------------

$ws.Range("B119").Select()
